$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 34: New York state hospitalization data for 16 April 2020.
$ws.Range("A34").Value = 43938
$ws.Range("A34").NumberFormat = "YYYY\-MM\-DD"

$ws.Range("B34").Value = -349
$ws.Range("C34").Value = -43
$ws.Range("D34").Value = -48
$ws.Range("F34").Value = 540
$ws.Range("G34").Value = 1915

# Match the author's cursor/selection state recorded in the diff.
$ws.Range("F35").Select()
